$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 97 reaction counts (Like +1, Comments +1)
$ws.Range("F97").Value = 5
$ws.Range("M97").Value = 4

# Append new row 98 - new Facebook post entry
$ws.Range("A98").Value = 40206
$ws.Range("A98").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B98").Value = 0.60347222222222219
$ws.Range("B98").NumberFormat = "h:mm:ss;@"
$ws.Range("C98").Value = "Public"
$ws.Range("D98").Value = "WE GOT SNOW 11 YEARS LATER"
$ws.Range("E98").Value = "10108028807642039"
$ws.Range("E98").NumberFormat = "@"
$ws.Range("F98").Value = 1
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 1
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0

# Restore frozen-pane view to the top of the sheet and move the selection
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$ws.Range("C12").Select()
